# Add new Mac-Address and Document Type rows to the
# master-reg_center_machine_device_h test data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data: regcntr_id, machine_id, device_id, lang_code, is_active, cr_by, cr_dtimes, eff_dtimes
$newRows = @(
    @(10002, 10032, 3000176, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000177, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000178, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000179, "eng", $true, "superadmin", "now()", "now()"),
    @(10002, 10032, 3000180, "eng", $true, "superadmin", "now()", "now()")
)

$startRow = 157
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $rowData[0]
    $ws.Cells.Item($r, 2).Value2 = $rowData[1]
    $ws.Cells.Item($r, 3).Value2 = $rowData[2]
    $ws.Cells.Item($r, 4).Value2 = $rowData[3]
    $ws.Cells.Item($r, 5).Value2 = $rowData[4]
    $ws.Cells.Item($r, 6).Value2 = $rowData[5]
    $ws.Cells.Item($r, 7).Value2 = $rowData[6]
    $ws.Cells.Item($r, 8).Value2 = $rowData[7]
}

# Switch calculation mode to manual, matching the workbook's calcPr change.
$excel.Calculation = -4135

# Update the active view/selection to match where the new data was entered.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 151
$ws.Range("D157").Select() | Out-Null
